# Weekly update: add a new "Coco" price record for Vega Modelo de Temuco.
# The new record is inserted as row 19 (most recent date first), pushing
# the previously existing rows 19-45 down to rows 20-46.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 19, shifting rows 19:45 down to 20:46.
$ws.Rows.Item(19).Insert()

# Populate the new row 19 with the new weekly record.
$ws.Cells.Item(19, 1).Value2  = 10
$ws.Cells.Item(19, 2).Value2  = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value2  = "La Araucanía"
$ws.Cells.Item(19, 4).Value2  = 44469
$ws.Cells.Item(19, 5).Value2  = 9
$ws.Cells.Item(19, 6).Value2  = "Fruta"
$ws.Cells.Item(19, 7).Value2  = 100108
$ws.Cells.Item(19, 8).Value2  = "Tropicales y subtropicales"
$ws.Cells.Item(19, 9).Value2  = 100108007
$ws.Cells.Item(19, 10).Value2 = "Coco"
$ws.Cells.Item(19, 11).Value2 = "Sin especificar"
$ws.Cells.Item(19, 12).Value2 = "Primera"
$ws.Cells.Item(19, 13).Value2 = 40
$ws.Cells.Item(19, 14).Value2 = 24000
$ws.Cells.Item(19, 15).Value2 = 24000
$ws.Cells.Item(19, 16).Value2 = 24000
$ws.Cells.Item(19, 17).Value2 = "$/malla 20 unidades"
$ws.Cells.Item(19, 18).Value2 = "Perú"
$ws.Cells.Item(19, 19).Value2 = 1200
$ws.Cells.Item(19, 20).Value2 = 20
